$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 16; $row++) {
    if ($row -eq 8) {
        $ws.Cells.Item($row, 5).Value = "FAIL"
    } else {
        $ws.Cells.Item($row, 5).Value = "PASS"
    }
}
